# Update cryptocurrency price/volume cells per the latest scrape run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.983.78"
$ws.Range("E2").Value = "  -2.29%  "
$ws.Range("D3").Value = "'2.098.87"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D5").Value = "'345.87"
$ws.Range("E5").Value = "  +2.35%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  -0.80%  "
$ws.Range("D7").Value = "'0.5156"
$ws.Range("E7").Value = "  -1.95%  "
$ws.Range("D8").Value = "'0.4411"
$ws.Range("E8").Value = "  -3.33%  "
$ws.Range("D9").Value = "'0.09237"
$ws.Range("E9").Value = "  +1.05%  "
$ws.Range("D10").Value = "'52.18"
$ws.Range("E10").Value = "  -5.15%  "
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").Value = "'25.22"
$ws.Range("E12").Value = "  +2.90%  "
$ws.Range("D13").Value = "'2.101.15"
$ws.Range("E13").Value = "  -0.69%  "
$ws.Range("D14").Value = "'6.734"
$ws.Range("E14").Value = "  -1.99%  "
$ws.Range("D15").Value = "'8.142"
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("D16").Value = "'99.25"
$ws.Range("E16").Value = "  +1.83%  "
$ws.Range("D17").Value = "'0.00001169"
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("D19").Value = "'20.65"
$ws.Range("E19").Value = "  +5.73%  "
$ws.Range("D20").Value = "'0.06662"
$ws.Range("E20").Value = "  -0.55%  "
$ws.Range("D22").Value = "'6.214"
$ws.Range("E22").Value = "  -1.86%  "
$ws.Range("D23").Value = "'30.078.64"
$ws.Range("E23").Value = "  -2.23%  "
$ws.Range("E24").Value = "  -2.45%  "
$ws.Range("D25").Value = "'2.338"
$ws.Range("E25").Value = "  -1.12%  "
$ws.Range("D26").Value = "'2.347.46"
$ws.Range("E26").Value = "  -0.70%  "
$ws.Range("D27").Value = "'21.95"
$ws.Range("E27").Value = "  -2.49%  "
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("D29").Value = "'162.33"
$ws.Range("E29").Value = "  -1.76%  "
$ws.Range("D30").Value = "'133.21"
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("D31").Value = "'1.166"
$ws.Range("E31").Value = "  -3.88%  "
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("D33").Value = "'1.633"
$ws.Range("E33").Value = "  -2.14%  "
$ws.Range("D34").Value = "'6.219"
$ws.Range("E34").Value = "  -2.69%  "
$ws.Range("D35").Value = "'3.957"
$ws.Range("D36").Value = "'6.177"
$ws.Range("E36").Value = "  +4.60%  "
$ws.Range("D37").Value = "'10.08"
$ws.Range("E37").Value = "  -5.48%  "
$ws.Range("D38").Value = "'0.02558"
$ws.Range("E38").Value = "  -3.57%  "
$ws.Range("D39").Value = "'0.06777"
$ws.Range("E39").Value = "  -1.61%  "
$ws.Range("D40").Value = "'0.2270"
$ws.Range("E40").Value = "  -2.51%  "
$ws.Range("D41").Value = "'12.49"
$ws.Range("E41").Value = "  -1.81%  "
$ws.Range("D42").Value = "'0.6894"
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("D43").Value = "'1.308"
$ws.Range("E43").Value = "  +3.73%  "
$ws.Range("D44").Value = "'0.6628"
$ws.Range("E44").Value = "  +1.83%  "
$ws.Range("D45").Value = "'14.13"
$ws.Range("E45").Value = "  -8.09%  "
$ws.Range("D46").Value = "'2.275"
$ws.Range("E46").Value = "  -2.31%  "
$ws.Range("D47").Value = "'3.628"
$ws.Range("E47").Value = "  -1.92%  "
$ws.Range("E48").Value = "  -5.29%  "
$ws.Range("D49").Value = "'1.219"
$ws.Range("E49").Value = "  -3.04%  "
$ws.Range("D50").Value = "'82.02"
$ws.Range("E50").Value = "  -1.92%  "
$ws.Range("D51").Value = "'0.07204"
$ws.Range("E51").Value = "  -1.54%  "
